$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 1051, pushing the existing data
# (rows 1051-1076) down to rows 1054-1079.
$ws.Rows("1051:1053").Insert()

# Fill in the 3 new rows with this week's price data for Tomate,
# "Larga vida" variety (Primera / Segunda / Tercera quality), keeping
# all the other columns identical to the surrounding rows.

# Row 1051: Larga vida - Primera
$ws.Range("A1051").Value = 2
$ws.Range("B1051").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C1051").Value = "Coquimbo"
$ws.Range("D1051").Value = 45239
$ws.Range("E1051").Value = 4
$ws.Range("F1051").Value = 100112020
$ws.Range("G1051").Value = "Tomate"
$ws.Range("H1051").Value = "Larga vida"
$ws.Range("I1051").Value = "Primera"
$ws.Range("J1051").Value = 2000
$ws.Range("K1051").Value = 14000
$ws.Range("L1051").Value = 15000
$ws.Range("M1051").Value = 14500
$ws.Range("N1051").Value = "$/bandeja 18 kilos"
$ws.Range("O1051").Value = "Provincia de Limarí"
$ws.Range("P1051").Value = 806
$ws.Range("Q1051").Value = 18
$ws.Range("R1051").Value = "Hortaliza"

# Row 1052: Larga vida - Segunda
$ws.Range("A1052").Value = 2
$ws.Range("B1052").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C1052").Value = "Coquimbo"
$ws.Range("D1052").Value = 45239
$ws.Range("E1052").Value = 4
$ws.Range("F1052").Value = 100112020
$ws.Range("G1052").Value = "Tomate"
$ws.Range("H1052").Value = "Larga vida"
$ws.Range("I1052").Value = "Segunda"
$ws.Range("J1052").Value = 1800
$ws.Range("K1052").Value = 11000
$ws.Range("L1052").Value = 12000
$ws.Range("M1052").Value = 11500
$ws.Range("N1052").Value = "$/bandeja 18 kilos"
$ws.Range("O1052").Value = "Provincia de Limarí"
$ws.Range("P1052").Value = 639
$ws.Range("Q1052").Value = 18
$ws.Range("R1052").Value = "Hortaliza"

# Row 1053: Larga vida - Tercera
$ws.Range("A1053").Value = 2
$ws.Range("B1053").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C1053").Value = "Coquimbo"
$ws.Range("D1053").Value = 45239
$ws.Range("E1053").Value = 4
$ws.Range("F1053").Value = 100112020
$ws.Range("G1053").Value = "Tomate"
$ws.Range("H1053").Value = "Larga vida"
$ws.Range("I1053").Value = "Tercera"
$ws.Range("J1053").Value = 1200
$ws.Range("K1053").Value = 7000
$ws.Range("L1053").Value = 8000
$ws.Range("M1053").Value = 7500
$ws.Range("N1053").Value = "$/bandeja 18 kilos"
$ws.Range("O1053").Value = "Provincia de Limarí"
$ws.Range("P1053").Value = 417
$ws.Range("Q1053").Value = 18
$ws.Range("R1053").Value = "Hortaliza"
